$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E) entirely, shifting
# reviews_average/latitude/longitude/is_permanently_closed/gmaps_link/
# latest_review_date left by one column (F:K -> E:J).
$ws.Range("E1").EntireColumn.Delete()
